$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15, pushing existing rows 15..46 down to 16..47
$ws.Rows.Item(15).Insert()

# Copy the formatting of the row above (row 14, a regular data row) onto the newly
# inserted row so it keeps the same look as the rest of the table (style index reused,
# rather than per-property copying which creates redundant style entries).
$ws.Range("A14:D14").Copy()
$ws.Range("A15:D15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the newly inserted row 15 with the new postulate entry
$ws.Cells.Item(15, 1).Value = "Z05_B01_P03"
$ws.Cells.Item(15, 2).Value = "Z05_B01"
$ws.Cells.Item(15, 3).Value = "Gleichstellung erreichen und alle Frauen und Mädchen zur Selbstbestimmung befähigen"
$ws.Cells.Item(15, 4).Value = "XXXGleichstellung erreichen und alle Frauen und Mädchen zur Selbstbestimmung befähigen"
